# Apply the "cards data" table-expansion edit to Лист1 (Sheet1) of the
# workbook: flips the header toggle cell A1 from "+" to "-" and appends
# five new product rows (6-10) with the same per-column formatting as the
# existing rows, then restores the original selection semantics.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Header toggle cell: "+" -> "-" -------------------------------------
$ws1.Range("A1").Value = "-"

# --- New rows 6-10 --------------------------------------------------------
# Column layout (same as rows 1-5): A text, B text, C text, D number,
# E text (numeric-looking), F text, G number (or quote-prefixed "--"),
# H number, I number, J text (plain, unformatted).

$rows = @(
    @{ Row=6;  A="-"; B="Somnambula";              C="jaws";     D=8; E="0.5";  F="Stout-Imperial`n Milk";      G=22;    H=20;   I=334; J="/src/img/png/product/somnambula.png" },
    @{ Row=7;  A="+"; B="Mission Of`nNutrition";    C="dieta";    D=5; E="0.5";  F="IPA-New England";            G="--";  H=12.5; I=307; J="/src/img/png/product/new_england.png" },
    @{ Row=8;  A="+"; B="nitro";                    C="jaws";     D=6; E="0.45"; F="Stout-Milk";                 G="--";  H=16;   I=249; J="/src/img/png/product/nitro.png" },
    @{ Row=9;  A="+"; B="stoner";                   C="Zagovor";  D=6; E="0.5";  F="Stout - Milk";               G=23;    H=17;   I=208; J="/src/img/png/product/stoner_milk.png" },
    @{ Row=10; A="+"; B="Tropical `nSlasher";       C="selfmade"; D=6; E="0.5";  F="Sour-Smoothie `n/Pastry";    G=20;    H=21;   I=368; J="/src/img/png/product/tropical_slasher.png" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Text columns: A, B, C, F (number format "@")
    foreach ($col in @("A", "B", "C", "F")) {
        $cell = $ws1.Range("$col$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r[$col]
    }

    # Numeric-looking columns stored as text: D, E, H, I use "0.00"
    # display format, but keep E's content ("0.5" / "0.45") as literal
    # text like the rest of the table (set as text, then switch the
    # number format so the stored value doesn't get coerced to a number).
    foreach ($col in @("D", "H", "I")) {
        $cell = $ws1.Range("$col$rowNum")
        $cell.NumberFormat = "0.00"
        $cell.Value = $r[$col]
    }

    $eCell = $ws1.Range("E$rowNum")
    $eCell.NumberFormat = "@"
    $eCell.Value = $r["E"]
    $eCell.NumberFormat = "0.00"

    # Column G: usually numeric, but two rows use the literal "--" text
    # (quote-prefixed, same pattern as G5 in the original table).
    $gCell = $ws1.Range("G$rowNum")
    $gCell.NumberFormat = "0.00"
    if ($r["G"] -eq "--") {
        $gCell.Value = "'--"
    } else {
        $gCell.Value = $r["G"]
    }

    # Column J: plain unformatted text (image path)
    $ws1.Range("J$rowNum").Value = $r["J"]

    # Typing multi-line text auto-expands the row height; the target
    # keeps the default row height, so re-fit it back down.
    $ws1.Rows.Item($rowNum).AutoFit()
}

# --- Selection: restore to the new active cell shown in the target ------
$ws1.Range("A6").Select()
